$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-5 with new TPM-derived values ---
# Row 2
$ws.Range("B2").Value = "Gdf7"
$ws.Range("C2").Value = "Bmpr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 0.02913733333333333
$ws.Range("H2").Value = 0.087412
$ws.Range("I2").Value = 0.2991461473965196
$ws.Range("J2").Value = 0.2991461473965196
$ws.Range("M2").Value = 46.29121633333333
$ws.Range("N2").Value = 138.873649
$ws.Range("O2").Value = 0.3133663986859022
$ws.Range("P2").Value = 0.3133663986859022
$ws.Range("Q2").Value = 1.348802600709778
$ws.Range("R2").Value = 12.139223406388
$ws.Range("S2").Value = 0.09374235089040941
$ws.Range("T2").Value = 0.09374235089040941

# Row 3
$ws.Range("B3").Value = "Gdf7"
$ws.Range("C3").Value = "Bmpr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("G3").Value = 0.02913733333333333
$ws.Range("H3").Value = 0.087412
$ws.Range("I3").Value = 0.2991461473965196
$ws.Range("J3").Value = 0.2991461473965196
$ws.Range("M3").Value = 46.81622333333333
$ws.Range("N3").Value = 140.44867
$ws.Range("O3").Value = 0.3169204109998198
$ws.Range("P3").Value = 0.3169204109998198
$ws.Range("Q3").Value = 1.364099904671111
$ws.Range("R3").Value = 12.27689914204
$ws.Range("S3").Value = 0.09480551998191765
$ws.Range("T3").Value = 0.09480551998191765

# Row 4
$ws.Range("B4").Value = "Gdf7"
$ws.Range("C4").Value = "Bmpr2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 0.02913733333333333
$ws.Range("H4").Value = 0.087412
$ws.Range("I4").Value = 0.2991461473965196
$ws.Range("J4").Value = 0.2991461473965196
$ws.Range("M4").Value = 38.53544233333333
$ws.Range("N4").Value = 115.606327
$ws.Range("O4").Value = 0.2608640200510233
$ws.Range("P4").Value = 0.2608640200510233
$ws.Range("Q4").Value = 1.122820028413778
$ws.Range("R4").Value = 10.105380255724
$ws.Range("S4").Value = 0.07803646659263205
$ws.Range("T4").Value = 0.07803646659263205

# Row 5
$ws.Range("B5").Value = "Gdf7"
$ws.Range("C5").Value = "Bmpr2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 0.02913733333333333
$ws.Range("H5").Value = 0.087412
$ws.Range("I5").Value = 0.2991461473965196
$ws.Range("J5").Value = 0.2991461473965196
$ws.Range("M5").Value = 16.07945366666667
$ws.Range("N5").Value = 48.238361
$ws.Range("O5").Value = 0.1088491702632547
$ws.Range("P5").Value = 0.1088491702632547
$ws.Range("Q5").Value = 0.4685124013035555
$ws.Range("R5").Value = 4.216611611732
$ws.Range("S5").Value = 0.03256180993156046
$ws.Range("T5").Value = 0.03256180993156046

# --- Add new rows 6-9 (Resolving-Mac as sending cluster) ---
# Row 6
$ws.Range("A6").Value = "Resolving-Mac"
$ws.Range("B6").Value = "Gdf7"
$ws.Range("C6").Value = "Bmpr2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.06826433333333333
$ws.Range("H6").Value = 0.204793
$ws.Range("I6").Value = 0.7008538526034804
$ws.Range("J6").Value = 0.7008538526034804
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 46.29121633333333
$ws.Range("N6").Value = 138.873649
$ws.Range("O6").Value = 0.3133663986859022
$ws.Range("P6").Value = 0.3133663986859022
$ws.Range("Q6").Value = 3.160039022184111
$ws.Range("R6").Value = 28.440351199657
$ws.Range("S6").Value = 0.2196240477954928
$ws.Range("T6").Value = 0.2196240477954928

# Row 7
$ws.Range("A7").Value = "Resolving-Mac"
$ws.Range("B7").Value = "Gdf7"
$ws.Range("C7").Value = "Bmpr2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.06826433333333333
$ws.Range("H7").Value = 0.204793
$ws.Range("I7").Value = 0.7008538526034804
$ws.Range("J7").Value = 0.7008538526034804
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 46.81622333333333
$ws.Range("N7").Value = 140.44867
$ws.Range("O7").Value = 0.3169204109998198
$ws.Range("P7").Value = 0.3169204109998198
$ws.Range("Q7").Value = 3.195878275034444
$ws.Range("R7").Value = 28.76290447531
$ws.Range("S7").Value = 0.2221148910179021
$ws.Range("T7").Value = 0.2221148910179021

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Gdf7"
$ws.Range("C8").Value = "Bmpr2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.06826433333333333
$ws.Range("H8").Value = 0.204793
$ws.Range("I8").Value = 0.7008538526034804
$ws.Range("J8").Value = 0.7008538526034804
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 38.53544233333333
$ws.Range("N8").Value = 115.606327
$ws.Range("O8").Value = 0.2608640200510233
$ws.Range("P8").Value = 0.2608640200510233
$ws.Range("Q8").Value = 2.630596280590111
$ws.Range("R8").Value = 23.675366525311
$ws.Range("S8").Value = 0.1828275534583913
$ws.Range("T8").Value = 0.1828275534583913

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Gdf7"
$ws.Range("C9").Value = "Bmpr2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.06826433333333333
$ws.Range("H9").Value = 0.204793
$ws.Range("I9").Value = 0.7008538526034804
$ws.Range("J9").Value = 0.7008538526034804
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 16.07945366666667
$ws.Range("N9").Value = 48.238361
$ws.Range("O9").Value = 0.1088491702632547
$ws.Range("P9").Value = 0.1088491702632547
$ws.Range("Q9").Value = 1.097653184919222
$ws.Range("R9").Value = 9.878878664273
$ws.Range("S9").Value = 0.07628736033169428
$ws.Range("T9").Value = 0.07628736033169428
